$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Switch to the Repayment Schedule sheet (was on Transactions)
$ws.Activate()

# Insert a new blank column before column N ("Late"), pushing Late/Heading/
# Outstanding one column to the right (N->O, O->P, P->Q)
$ws.Columns.Item(14).Insert()

# Match column N's width to the width inherited from its left neighbour (M)
$ws.Columns.Item(14).ColumnWidth = 10.33

# Leave the selection on the newly relevant column near the bottom of the data
$ws.Range("P10").Select()
